# Update the "想去人数" (F column) counts on several sheets to reflect the
# refreshed data snapshot published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet "展览"
# -----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 3414
$ws1.Range("F4").Value  = 580
$ws1.Range("F10").Value = 634
$ws1.Range("F12").Value = 426
$ws1.Range("F13").Value = 67
$ws1.Range("F14").Value = 492
$ws1.Range("F15").Value = 328
$ws1.Range("F16").Value = 59
$ws1.Range("F19").Value = 184

# -----------------------------------------------------------------------
# Sheet "本地生活"
# -----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6205
$ws3.Range("F5").Value = 1781

# -----------------------------------------------------------------------
# Sheet "全部类型"
# -----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6205
$ws4.Range("F5").Value  = 1781
$ws4.Range("F6").Value  = 3414
$ws4.Range("F9").Value  = 580
$ws4.Range("F20").Value = 634
$ws4.Range("F24").Value = 426
$ws4.Range("F26").Value = 67
$ws4.Range("F27").Value = 492
$ws4.Range("F29").Value = 328
$ws4.Range("F30").Value = 59
$ws4.Range("F40").Value = 184
